# Update to correspond to performance changes in pymassspec and clean up
#
# The UserLibrary sheet had a batch of columns (RI non-polar/RI polar, and
# the whole quality/odor/synonyms/Reference/Conversion Notes/DB-5/FFAP/
# bLCP/aLCP/Synonyms/Comments block) that are no longer produced by the
# updated pymassspec pipeline, so they are cleared out. A new "Formula"
# column (O) is now also populated for the first couple of library rows,
# duplicating the Nickname value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-unused header columns + their (empty) data underneath.
# Using ClearContents (not Delete) so cells to the right are NOT shifted.
$ws.Range("V1:W1").ClearContents()
$ws.Range("Z1:AD1").ClearContents()
$ws.Range("AG1:AM1").ClearContents()

# Populate the "Formula" column (O) for the two rows that now carry it,
# mirroring the Nickname column (H).
$ws.Range("O2").Value = $ws.Range("H2").Value()
$ws.Range("O3").Value = $ws.Range("H3").Value()

# Matches the saved selection in the edited workbook.
[void]$ws.Range("O4").Select()
